$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix Lab 8 testbench bugs: negate the High-Pass coefficient values in column D
# for the symmetric taps that were incorrectly entered as positive.
$ws.Range("D4").Value = -0.0031
$ws.Range("D5").Value = -0.0108
$ws.Range("D9").Value = -0.0807
$ws.Range("D10").Value = -0.2913
$ws.Range("D12").Value = -0.2913
$ws.Range("D13").Value = -0.0807
$ws.Range("D17").Value = -0.0108
$ws.Range("D18").Value = -0.0031

# Update the view: scroll so row 4 is at the top, zoom to 190%, and select E12
$ws.Range("E12").Select()
$excel.ActiveWindow.Zoom = 190
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
